# Avance muestreo de datos SpaOnline
# Populate the "Servicio" sample-data sheet: drop the unused "Fabricante"
# column, rename the "Nombre" header to "Nombre Servicio", and fill in the
# three sample service rows (with the "Combinacion unica" helper formula
# following the name column left after Fabricante is removed).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Servicio")
$ws.Activate()

# Remove the "Fabricante" column (old column C) entirely - this shifts the
# "Combinacion unica" column from D into C and drops the stale refs.
$ws.Columns.Item(3).Delete()

# Header row
$ws.Range("B1").Value = "Nombre Servicio"

# Sample data rows
$ws.Range("B2").Value = "Limpieza facial"
$ws.Range("C2").Formula = "=+B2"

$ws.Range("B3").Value = "Masaje completo"
$ws.Range("B4").Value = "Chocolaterapia"
$ws.Range("C3:C4").Formula = "=+B3"

$ws.Range("B2").Select()
